$d = $word.ActiveDocument

# 1) Remove the "Counsel for the State..." narrative sentence (and trailing
#    spaces) that was inserted after "Defendant appeared in Court for
#    arraignment on December 12, 2021."  The text spans what used to be four
#    separate runs; Find/Replace across the range collapses them.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Counsel for the State of Ohio made a motion to amend the charge of DUS UCM to Assured Clear Distrance Ahead. The Court found the amendment did alter the name or identify of the offense and therefore the motion is Denied.   ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2) The offense table cell no longer reflects an amendment.
$find.Execute("DUS UCM - AMENDED", $true, $false, $false, $false, $false,
              $true, 1, $false, "DUS UCM", 2)

# 3) The plea for each count changed from "No Contest" to "Guilty".
$find.Execute("No Contest", $true, $false, $false, $false, $false,
              $true, 1, $false, "Guilty", 2)
